# Add a new "Save" column (H) to the s_vals sheet, mirroring the
# formatting of the existing header cells (e.g. G1) and appending a
# matching data value of 0 for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save"
$ws.Range("H1").Value = "Save"

# Copy the formatting (bold font + border, centered/top-aligned) from the
# neighboring header cell G1 onto H1 so it matches the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell H2 = 0
$ws.Range("H2").Value = 0
